$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D2").Value = 35
$wsForecast.Range("D7").Value = 37

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "1072"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "491"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "233"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "35"
